# Update countries & provincias Spain
# Applies the 3-Jul-2020 10:50 data refresh to the "Pais" sheet:
#  - refreshed case/death/recovery counters for several countries
#  - Namibia's case count overtook Angola/Siria/Birmania/Comoras, shifting
#    those four rows down by one position
#  - Laos/Santa Lucia and Dominica/Fiyi swapped order (tied counts)
#  - updated "datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 10:50"
$ws.Range("B6").Value = 667883
$ws.Range("C6").Value = 6718
$ws.Range("D6").Value = 437893
$ws.Range("E6").Value = 220131
$ws.Range("G6").Value = 176
$ws.Range("H6").Value = 9859
$ws.Range("B7").Value = 628205
$ws.Range("C7").Value = 1037
$ws.Range("D7").Value = 380374
$ws.Range("E7").Value = 229590
$ws.Range("G7").Value = 16
$ws.Range("H7").Value = 18241
$ws.Range("D45").Value = 22651
$ws.Range("E45").Value = 11003
$ws.Range("B47").Value = 32324
$ws.Range("C47").Value = 302
$ws.Range("D47").Value = 17331
$ws.Range("E47").Value = 14174
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 819
$ws.Range("B51").Value = 27611
$ws.Range("C51").Value = 564
$ws.Range("D51").Value = 17607
$ws.Range("E51").Value = 9678
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 326
$ws.Range("B100").Value = 3148
$ws.Range("C100").Value = 68
$ws.Range("D100").Value = 463
$ws.Range("E100").Value = 2675
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 10
$ws.Range("B116").Value = 1828
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 1539
$ws.Range("E116").Value = 210
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 79
$ws.Range("B118").Value = 1720
$ws.Range("C118").Value = 20
$ws.Range("E118").Value = 226
$ws.Range("B120").Value = 1650
$ws.Range("C120").Value = 16
$ws.Range("E120").Value = 155
$ws.Range("D154").Value = 253
$ws.Range("E154").Value = 281
$ws.Range("A161").Value = "Namibia"
$ws.Range("B161").Value = 335
$ws.Range("C161").Value = 42
$ws.Range("D161").Value = 25
$ws.Range("E161").Value = 310
$ws.Range("H161").Value = 0
$ws.Range("A162").Value = "Angola"
$ws.Range("B162").Value = 315
$ws.Range("D162").Value = 97
$ws.Range("E162").Value = 201
$ws.Range("H162").Value = 17
$ws.Range("A163").Value = "Siria"
$ws.Range("B163").Value = 312
$ws.Range("D163").Value = 113
$ws.Range("E163").Value = 190
$ws.Range("H163").Value = 9
$ws.Range("A164").Value = "Birmania"
$ws.Range("B164").Value = 304
$ws.Range("D164").Value = 223
$ws.Range("E164").Value = 75
$ws.Range("H164").Value = 6
$ws.Range("A165").Value = "Comoras"
$ws.Range("B165").Value = 303
$ws.Range("D165").Value = 200
$ws.Range("E165").Value = 96
$ws.Range("H165").Value = 7
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"
